$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# The "createdste" column (E2:E6) previously held a different date string
# per row (20-02-2022 .. 20-06-2022). The bulk-upload fix now stamps every
# row with the same date, "2022-10-20", stored as text (column E uses a
# text number format already).
$ws.Range("E2:E6").NumberFormat = "@"
$ws.Range("E2").Value = "2022-10-20"
$ws.Range("E3").Value = "2022-10-20"
$ws.Range("E4").Value = "2022-10-20"
$ws.Range("E5").Value = "2022-10-20"
$ws.Range("E6").Value = "2022-10-20"

# Update the saved cursor/selection position recorded in the sheet view.
$ws.Activate()
$ws.Range("G6").Select()
